$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 252.5340982325772

$ws.Range("B3").Value = 0.006876353814593728
$ws.Range("C3").Value = 86.29678392075563
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 92.93799413385902

$ws.Range("B4").Value = 0.3464964993005633
$ws.Range("C4").Value = 0.05231270169004087
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 246.9852506941017
$ws.Range("G4").Value = 250.4666593217958

$ws.Range("B5").Value = 3.182878228561681
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 5.488907176552729

$ws.Range("B6").Value = 1.505614041169197
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 3.811642989160245
